$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slightly-updated timestamp on the existing last row (row 18)
$ws.Range("A18").Value = 45865.79197570602

# Append the new row captured by the scheduled task
$ws.Range("A19").Value = 45865.83358738009
$ws.Range("A19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B19").Value = 2025
$ws.Range("C19").Value = 30
$ws.Range("D19").Value = 14.22
$ws.Range("E19").Value = 89.29000000000001
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 4.42
$ws.Range("H19").Value = "ESE"
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = "20:00:21"
